# Update column G ("K") values on Sheet1, rows 2-7, to reflect the
# regenerated save_data using K (strike count) instead of Strike#.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 9
$ws.Range("G5").Value = 5
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 2
